$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Widget E
$ws.Range("A2").Value = "82457957-3d05-4671-99dc-975f83a8a276"
$ws.Range("B2").Value = "Widget E"
$ws.Range("C2").Value = "A heavily used widget"
$ws.Range("D2").Value = 25
$ws.Range("E2").Value = 10

# Row 3: Widget A
$ws.Range("A3").Value = "4d3bfcc6-4ba2-40be-9819-27ab3369c59f"
$ws.Range("B3").Value = "Widget A"
$ws.Range("C3").Value = "A Premium Widget"
$ws.Range("D3").Value = 10
$ws.Range("E3").Value = 15

# Row 4: Dropdown D
$ws.Range("A4").Value = "8abd302b-3202-41d5-9cdf-d308f77e7c95"
$ws.Range("B4").Value = "Dropdown D"
$ws.Range("C4").Value = "A premium dropdown"
$ws.Range("D4").Value = 50
$ws.Range("E4").Value = 5
